$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update payout values (Bet column) for several rows
$ws.Range("B14").Value = 100
$ws.Range("B16").Value = 4
$ws.Range("B18").Value = 20
$ws.Range("B20").Value = 40
$ws.Range("B21").Value = 45

# Apply the "1:" custom number format to the Chance column (C14:C25)
$ws.Range("C14:C25").NumberFormat = """1:""#"

# Update the selection to F19
$ws.Range("F19").Select()
